# Adds columns I (I0) and J (IF) to Sheet1, mirroring the structure of the
# existing H (IP) column: header cell styled like the other header cells,
# and plain numeric values for each data row (rows 2-48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1) - same text + formatting style as the other headers (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Data values for column I (I0), rows 2-48
$iValues = @(8,5,7,5,7,8,7,2,2,5,6,10,5,6,8,8,8,4,5,7,6,8,7,7,7,7,5,6,9,8,7,7,6,7,3,8,7,6,8,8,1,1,1,2,1,1,4)

# Data values for column J (IF), rows 2-48
$jValues = @(8,6,7,5,7,8,8,3,4,5,6,11,6,7,8,8,9,4,5,7,6,8,7,7,7,7,5,6,10,8,9,8,6,7,5,8,8,6,8,8,1,4,6,6,3,3,4)

for ($r = 2; $r -le 48; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
